$d = $word.ActiveDocument

# Change 1: "Την με αριθ. Φ.353.1" -> "Τη με αριθ. Φ.353.1"
$d.Content.Find.Execute("Την με αριθ. Φ.353.1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Τη με αριθ. Φ.353.1", 2)

# Change 2: replace old Φ.350.2 reference (across 3 runs) with new Φ.351.1 reference (single run)
$d.Content.Find.Execute("Την με αριθ. Φ.350.2/1/32958/Ε3/27-2-2018  (ΑΔΑ:6Π414653ΠΣ-7ΕΝ) Υπουργική Απόφαση με θέμα: «Τοποθέτηση Περιφερειακών Διευθυντών Εκπαίδευσης»", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Τη με αριθ. Φ.351.1/11/48020/Ε3/28-3-2019 (ΑΔΑ: ΩΩΤΗ4653ΠΣ-ΒΔ3) Υπουργική Απόφαση με θέμα: «Τοποθέτηση Περιφερειακών Διευθυντών Εκπαίδευσης»", 2)
